# Commit: "Changed the 'goto' statement on the EXCEL template"
#
# On Sheet1, row 4 ("What is the patient's sex?") the Then_Goto / Else_Goto
# columns (I4 / J4) used to hold the literal target variable names
# ("pregnant" / "onset_date") as text. The template now points to the
# destination row numbers instead, so I4/J4 become plain numbers (4 / 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5

$ws.Range("J6").Select()
